$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3..21: the date text in column A switches from dd/mm/yyyy to dd-mm-yyyy.
# A handful of rows also get updated attendance-count values (D/E/G/H).
$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; G = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 7;  Date = "11-08-2022" },
    @{ Row = 8;  Date = "15-08-2022" },
    @{ Row = 9;  Date = "18-08-2022" },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 11; Date = "25-08-2022" },
    @{ Row = 12; Date = "29-08-2022" },
    @{ Row = 13; Date = "01-09-2022" },
    @{ Row = 14; Date = "05-09-2022" },
    @{ Row = 15; Date = "08-09-2022" },
    @{ Row = 16; Date = "12-09-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 17; Date = "15-09-2022" },
    @{ Row = 18; Date = "19-09-2022" },
    @{ Row = 19; Date = "22-09-2022" },
    @{ Row = 20; Date = "26-09-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 21; Date = "29-09-2022" }
)

foreach ($r in $rows) {
    # Force the cell to text format before writing so Excel doesn't
    # auto-convert the dd-mm-yyyy-looking string into a date serial,
    # then strip the format again so the cell keeps its original
    # (default) style, matching the un-styled inline-string cell.
    $cell = $ws.Range("A$($r.Row)")
    $cell.NumberFormat = "@"
    $cell.Value = $r.Date
    $cell.ClearFormats()

    if ($r.ContainsKey("D")) { $ws.Range("D$($r.Row)").Value = $r.D }
    if ($r.ContainsKey("E")) { $ws.Range("E$($r.Row)").Value = $r.E }
    if ($r.ContainsKey("G")) { $ws.Range("G$($r.Row)").Value = $r.G }
    if ($r.ContainsKey("H")) { $ws.Range("H$($r.Row)").Value = $r.H }
}
